$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap values in columns A, B, D, E, F, G, H between row 3 and row 4
$cols = @("A", "B", "D", "E", "F", "G", "H")

foreach ($col in $cols) {
    $addr3 = "$col`3"
    $addr4 = "$col`4"
    $val3 = $ws.Range($addr3).Value2
    $val4 = $ws.Range($addr4).Value2
    $ws.Range($addr3).Value2 = $val4
    $ws.Range($addr4).Value2 = $val3
}
